$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the casing of the "Password" header label to lowercase "password"
$ws.Range("B1").Value = "password"

# Remove the bold styling from the header row (A1:B1)
$ws.Range("A1:B1").Font.Bold = $false

# Leave the cursor/selection on B1
$ws.Range("B1").Select()
